# BOT; UPDATE DATA
# Adds the next day's row (5/4/2020) to the "相談件数" (consultation count)
# sheet, pushing the trailing footnote row down by one, and extends the
# print area / dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Push the footnote that currently lives in row 100 down to row 101 ---
# Row 100 today only has a single cell (B100, the "※4/8..." footnote).
# Copy it (value + number format) down to B101 first so we don't clobber it
# when row 100 becomes a real data row below.
$ws.Range("B100").Copy($ws.Range("B101"))

# Give A101 the same date-style formatting as the rest of column A
# (an empty, but styled, cell) by cloning A99's format then clearing the
# value that the copy brought along.
$ws.Range("A99").Copy($ws.Range("A101"))
$ws.Range("A101").ClearContents()

# --- 2. Turn row 100 into a real data row for 2020-05-04 ---
# Seed it from row 99 so it inherits the same per-column number formats
# (date / "0_);[Red](0)" / right-aligned wrapped general), then overwrite
# with the new day's values.
$ws.Range("A99:E99").Copy($ws.Range("A100:E100"))
$ws.Range("A100").Value = 43955
$ws.Range("B100").Value = 428
$ws.Range("C100").Value = 33464
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = 6958

# --- 3. Update the print area to include the new last row ---
$printArea = $wb.Names.Item($ws.Name + "!Print_Area")
$printArea.RefersTo = '=相談件数!$A$1:$E$102'

# --- 4. Move the active selection to the new last cell, matching the
#        author's saved view state ---
$ws.Range("E101").Select() | Out-Null
